$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.ClearContents()

# Header row (row 1) - values only, style already present on cells
$ws.Cells.Item(1, 1).Value = 'Sending cluster'
$ws.Cells.Item(1, 2).Value = 'Ligand symbol'
$ws.Cells.Item(1, 3).Value = 'Receptor symbol'
$ws.Cells.Item(1, 4).Value = 'Target cluster'
$ws.Cells.Item(1, 5).Value = 'Ligand-expressing cells'
$ws.Cells.Item(1, 6).Value = 'Ligand detection rate'
$ws.Cells.Item(1, 7).Value = 'Ligand average expression value'
$ws.Cells.Item(1, 8).Value = 'Ligand total expression value'
$ws.Cells.Item(1, 9).Value = 'Ligand derived specificity of average expression value'
$ws.Cells.Item(1, 10).Value = 'Ligand derived specificity of total expression value'
$ws.Cells.Item(1, 11).Value = 'Receptor-expressing cells'
$ws.Cells.Item(1, 12).Value = 'Receptor detection rate'
$ws.Cells.Item(1, 13).Value = 'Receptor average expression value'
$ws.Cells.Item(1, 14).Value = 'Receptor total expression value'
$ws.Cells.Item(1, 15).Value = 'Receptor derived specificity of average expression value'
$ws.Cells.Item(1, 16).Value = 'Receptor derived specificity of total expression value'
$ws.Cells.Item(1, 17).Value = 'Edge average expression weight'
$ws.Cells.Item(1, 18).Value = 'Edge total expression weight'
$ws.Cells.Item(1, 19).Value = 'Edge average expression derived specificity'
$ws.Cells.Item(1, 20).Value = 'Edge total expression derived specificity'

# Data rows 2-7
# Row 2
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Ccl4'
$ws.Cells.Item(2, 3).Value = 'Ccr3'
$ws.Cells.Item(2, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.2400906666666667
$ws.Cells.Item(2, 8).Value = 0.720272
$ws.Cells.Item(2, 9).Value = 0.001631540293869566
$ws.Cells.Item(2, 10).Value = 0.001631540293869566
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.08962966666666666
$ws.Cells.Item(2, 14).Value = 0.268889
$ws.Cells.Item(2, 15).Value = 0.4339761198462219
$ws.Cells.Item(2, 16).Value = 0.4339761198462219
$ws.Cells.Item(2, 17).Value = 0.02151924642311111
$ws.Cells.Item(2, 18).Value = 0.193673217808
$ws.Cells.Item(2, 19).Value = 0.000708049526106279
$ws.Cells.Item(2, 20).Value = 0.0007080495261062789

# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Ccl4'
$ws.Cells.Item(3, 3).Value = 'Ccr3'
$ws.Cells.Item(3, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.2400906666666667
$ws.Cells.Item(3, 8).Value = 0.720272
$ws.Cells.Item(3, 9).Value = 0.001631540293869566
$ws.Cells.Item(3, 10).Value = 0.001631540293869566
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.1169016666666667
$ws.Cells.Item(3, 14).Value = 0.350705
$ws.Cells.Item(3, 15).Value = 0.5660238801537781
$ws.Cells.Item(3, 16).Value = 0.5660238801537781
$ws.Cells.Item(3, 17).Value = 0.02806699908444445
$ws.Cells.Item(3, 18).Value = 0.25260299176
$ws.Cells.Item(3, 19).Value = 0.0009234907677632874
$ws.Cells.Item(3, 20).Value = 0.0009234907677632873

# Row 4
$ws.Cells.Item(4, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4, 2).Value = 'Ccl4'
$ws.Cells.Item(4, 3).Value = 'Ccr3'
$ws.Cells.Item(4, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 91.82408133333333
$ws.Cells.Item(4, 8).Value = 275.472244
$ws.Cells.Item(4, 9).Value = 0.623992138981758
$ws.Cells.Item(4, 10).Value = 0.623992138981758
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.08962966666666666
$ws.Cells.Item(4, 14).Value = 0.268889
$ws.Cells.Item(4, 15).Value = 0.4339761198462219
$ws.Cells.Item(4, 16).Value = 0.4339761198462219
$ws.Cells.Item(4, 17).Value = 8.230161801879554
$ws.Cells.Item(4, 18).Value = 74.07145621691599
$ws.Cells.Item(4, 19).Value = 0.2707976872898477
$ws.Cells.Item(4, 20).Value = 0.2707976872898477

# Row 5
$ws.Cells.Item(5, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(5, 2).Value = 'Ccl4'
$ws.Cells.Item(5, 3).Value = 'Ccr3'
$ws.Cells.Item(5, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 91.82408133333333
$ws.Cells.Item(5, 8).Value = 275.472244
$ws.Cells.Item(5, 9).Value = 0.623992138981758
$ws.Cells.Item(5, 10).Value = 0.623992138981758
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1169016666666667
$ws.Cells.Item(5, 14).Value = 0.350705
$ws.Cells.Item(5, 15).Value = 0.5660238801537781
$ws.Cells.Item(5, 16).Value = 0.5660238801537781
$ws.Cells.Item(5, 17).Value = 10.73438814800222
$ws.Cells.Item(5, 18).Value = 96.60949333201999
$ws.Cells.Item(5, 19).Value = 0.3531944516919103
$ws.Cells.Item(5, 20).Value = 0.3531944516919103

# Row 6
$ws.Cells.Item(6, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(6, 2).Value = 'Ccl4'
$ws.Cells.Item(6, 3).Value = 'Ccr3'
$ws.Cells.Item(6, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 55.09165833333333
$ws.Cells.Item(6, 8).Value = 165.274975
$ws.Cells.Item(6, 9).Value = 0.3743763207243725
$ws.Cells.Item(6, 10).Value = 0.3743763207243724
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.08962966666666666
$ws.Cells.Item(6, 14).Value = 0.268889
$ws.Cells.Item(6, 15).Value = 0.4339761198462219
$ws.Cells.Item(6, 16).Value = 0.4339761198462219
$ws.Cells.Item(6, 17).Value = 4.937846972530555
$ws.Cells.Item(6, 18).Value = 44.440622752775
$ws.Cells.Item(6, 19).Value = 0.1624703830302679
$ws.Cells.Item(6, 20).Value = 0.1624703830302678

# Row 7
$ws.Cells.Item(7, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(7, 2).Value = 'Ccl4'
$ws.Cells.Item(7, 3).Value = 'Ccr3'
$ws.Cells.Item(7, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 55.09165833333333
$ws.Cells.Item(7, 8).Value = 165.274975
$ws.Cells.Item(7, 9).Value = 0.3743763207243725
$ws.Cells.Item(7, 10).Value = 0.3743763207243724
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.1169016666666667
$ws.Cells.Item(7, 14).Value = 0.350705
$ws.Cells.Item(7, 15).Value = 0.5660238801537781
$ws.Cells.Item(7, 16).Value = 0.5660238801537781
$ws.Cells.Item(7, 17).Value = 6.440306678597223
$ws.Cells.Item(7, 18).Value = 57.962760107375
$ws.Cells.Item(7, 19).Value = 0.2119059376941046
$ws.Cells.Item(7, 20).Value = 0.2119059376941046
